$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9 = @{ B = 43; C = 88; D = 48.86363636363637 }
    10 = @{ B = 52; C = 91; D = 57.14285714285714 }
    11 = @{ B = 63; C = 122; D = 51.63934426229508 }
    12 = @{ B = 61; C = 100; D = 61 }
    13 = @{ B = 64; C = 125; D = 51.2 }
    14 = @{ B = 70; C = 128; D = 54.6875 }
    15 = @{ B = 79; C = 139; D = 56.83453237410072 }
    16 = @{ B = 78; C = 127; D = 61.41732283464567 }
    17 = @{ B = 79; C = 117; D = 67.52136752136752 }
    18 = @{ B = 85; C = 143; D = 59.44055944055944 }
    19 = @{ B = 91; C = 143; D = 63.63636363636363 }
    20 = @{ B = 89; C = 146; D = 60.95890410958904 }
    21 = @{ B = 77; C = 139; D = 55.39568345323741 }
    22 = @{ B = 88; C = 112; D = 78.57142857142857 }
    24 = @{ B = 71; C = 104; D = 68.26923076923077 }
    25 = @{ B = 76; C = 108; D = 70.37037037037037 }
    26 = @{ B = 69; C = 127; D = 54.33070866141733 }
    27 = @{ B = 74; C = 98; D = 75.51020408163265 }
    28 = @{ B = 75; C = 123; D = 60.97560975609756 }
    31 = @{ B = 80; C = 121; D = 66.11570247933885 }
    33 = @{ B = 107; C = 141; D = 75.88652482269504 }
    34 = @{ B = 103; C = 143; D = 72.02797202797203 }
    35 = @{ B = 101; C = 138; D = 73.18840579710145 }
    36 = @{ B = 108; C = 143; D = 75.52447552447552 }
    37 = @{ B = 115; C = 160; D = 71.875 }
    38 = @{ B = 115; C = 156; D = 73.71794871794873 }
    39 = @{ B = 110; C = 164; D = 67.07317073170732 }
    41 = @{ B = 105; C = 143; D = 73.42657342657343 }
    42 = @{ B = 123; C = 182; D = 67.58241758241759 }
    43 = @{ B = 121; C = 161; D = 75.15527950310559 }
    44 = @{ B = 136; C = 189; D = 71.95767195767195 }
    45 = @{ B = 143; C = 197; D = 72.58883248730965 }
    46 = @{ B = 141; C = 209; D = 67.46411483253588 }
    47 = @{ B = 150; C = 215; D = 69.76744186046511 }
    48 = @{ B = 165; C = 225; D = 73.33333333333333 }
    49 = @{ B = 172; C = 226; D = 76.10619469026548 }
    50 = @{ B = 169; C = 227; D = 74.44933920704845 }
    51 = @{ B = 151; C = 242; D = 62.39669421487604 }
    52 = @{ B = 123; C = 197; D = 62.43654822335025 }
    53 = @{ B = 145; C = 195; D = 74.35897435897436 }
    54 = @{ B = 159; C = 216; D = 73.61111111111111 }
    55 = @{ B = 143; C = 232; D = 61.63793103448276 }
    56 = @{ B = 138; C = 203; D = 67.98029556650246 }
    57 = @{ B = 159; C = 222; D = 71.62162162162163 }
    58 = @{ B = 64; C = 244; D = 26.22950819672131 }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $u.B
    $ws.Cells.Item($row, 3).Value = $u.C
    $ws.Cells.Item($row, 4).Value = $u.D
}
